$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
# Price cells (column D) are stored as plain text in the workbook, even when
# the text looks like a number (e.g. "4.88", "1.00", "0.0924"). A leading
# apostrophe forces Excel to keep the entered value as text (quote-prefixed)
# instead of silently coercing it to a numeric cell.
$ws.Range("D2").Value = '69.507.40'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '2.493.92'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''569.32'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '''165.98'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("D12").Value = '''4.88'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '2.950.08'
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("D14").Value = '69.555.54'
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").Value = '''0.0000175'
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("D16").Value = '''24.15'
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").Value = '2.483.06'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  -1.14%  '
$ws.Range("D20").Value = '''352.60'
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("E22").Value = '  -3.51%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("E26").Value = '  -1.25%  '
$ws.Range("E27").Value = '  -2.39%  '
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("D31").Value = '''3.55'
$ws.Range("E31").Value = '  +139.47%  '
$ws.Range("E32").Value = '  -3.07%  '
$ws.Range("D33").Value = '''439.08'
$ws.Range("E33").Value = '  -4.77%  '
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("E36").Value = '  -3.20%  '
$ws.Range("D37").Value = '''152.93'
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("D39").Value = '''18.14'
$ws.Range("E39").Value = '  -1.60%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("D42").Value = '''4.59'
$ws.Range("E42").Value = '  -2.29%  '
$ws.Range("E43").Value = '  -1.84%  '
$ws.Range("E44").Value = '  -2.19%  '
$ws.Range("E45").Value = '  -3.81%  '
$ws.Range("D46").Value = '''139.09'
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("E48").Value = '  -2.54%  '
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("E50").Value = '  -0.79%  '
$ws.Range("D51").Value = '''0.0924'
$ws.Range("E51").Value = '  -0.44%  '
